$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Event:" dropdown value in B3 (merged B3:C3) from "Auto Renewal Letter" to "Sort Renewal List"
$ws.Range("B3").Value = "Sort Renewal List"

# Move the active selection to D4 (matches post-edit cursor position)
$ws.Range("D4").Select()
